# Add 2022-Q1 fund-holdings data, and refresh the "总计" (totals) summary
# sheet with the new quarter's row.
#
# The existing "总计" sheet (last tab) is repurposed in place to become the
# new "2022-Q1" data sheet (keeping its sheetId/r:id, just renamed and
# re-populated), and a brand-new "总计" sheet is appended right after it
# with the refreshed totals table - this mirrors how the sheetId/r:id
# numbering shifts in the target workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Repurpose the current last sheet ("总计") into the "2022-Q1" sheet.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(5)          # "2021-Q4" - same table layout
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Pull the header-row and index-column formatting (bold, centered, thin
# border) from an existing fund-holdings sheet so the new sheet matches it
# exactly (same style slot).
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)        # xlPasteFormats

$template.Range("A2").Copy()
$q1.Range("A2:A23").PasteSpecial(-4122)       # xlPasteFormats

$fundHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q1.Cells.Item(1, $col).Value = $fundHeaders[$col - 2]
}

# Columns B (fund code) and D:G (scale/position figures) are stored as text
# in the source data (leading zeros in fund codes must survive), so format
# them as Text before writing the values.
$q1.Range("B2:B23").NumberFormat = "@"
$q1.Range("D2:G23").NumberFormat = "@"

$fundRows = @(
    @("040008", "华安策略优选混合", "52.77", "92.07", "6.26", "3.3034", 4),
    @("010787", "华安优势企业混合A", "36.73", "92.23", "4.13", "1.5169", 8),
    @("519690", "交银稳健配置混合A", "16.44", "88.28", "3.15", "0.5179", 6),
    @("960017", "交银稳健配置混合H", "16.44", "88.28", "3.15", "0.5179", 6),
    @("515210", "国泰中证钢铁ETF", "16.24", "99.25", "2.77", "0.4498", 10),
    @("502023", "鹏华国证钢铁行业指数（LOF）", "15.55", "94.76", "2.44", "0.3794", 10),
    @("005273", "华商可转债债券A", "10.06", "39.60", "2.85", "0.2867", 2),
    @("011238", "华安聚恒精选混合A", "5.14", "91.92", "4.13", "0.2123", 8),
    @("005521", "华安红利精选混合", "4.72", "91.71", "4.14", "0.1954", 10),
    @("010296", "万家互联互通中国优势量化策略混合A", "5.55", "92.06", "3.30", "0.1832", 4),
    @("005284", "华商可转债债券C", "3.80", "39.60", "2.85", "0.1083", 2),
    @("168203", "中融国证钢铁行业指数", "4.30", "92.58", "2.37", "0.1019", 10),
    @("010788", "华安优势企业混合C", "2.23", "92.23", "4.13", "0.0921", 8),
    @("002434", "中银宏利灵活配置混合A", "5.86", "30.78", "0.71", "0.0416", 5),
    @("003966", "中银润利灵活配置混合A", "5.69", "25.96", "0.63", "0.0358", 4),
    @("002261", "中银宝利灵活配置混合A", "4.35", "31.13", "0.71", "0.0309", 5),
    @("011239", "华安聚恒精选混合C", "0.65", "91.92", "4.13", "0.0268", 8),
    @("003967", "中银润利灵活配置混合C", "3.85", "25.96", "0.63", "0.0243", 4),
    @("010297", "万家互联互通中国优势量化策略混合C", "0.53", "92.06", "3.30", "0.0175", 4),
    @("002435", "中银宏利灵活配置混合C", "2.33", "30.78", "0.71", "0.0165", 5),
    @("002262", "中银宝利灵活配置混合C", "1.60", "31.13", "0.71", "0.0114", 5),
    @("515510", "嘉实中证500成长估值ETF", "0.15", "98.79", "1.23", "0.0018", 7)
)

$row = 2
foreach ($fund in $fundRows) {
    $q1.Cells.Item($row, 1).Value = $row - 2
    $q1.Cells.Item($row, 2).Value = $fund[0]
    $q1.Cells.Item($row, 3).Value = $fund[1]
    $q1.Cells.Item($row, 4).Value = $fund[2]
    $q1.Cells.Item($row, 5).Value = $fund[3]
    $q1.Cells.Item($row, 6).Value = $fund[4]
    $q1.Cells.Item($row, 7).Value = $fund[5]
    $q1.Cells.Item($row, 8).Value = $fund[6]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2) Append a fresh "总计" sheet after "2022-Q1" with the refreshed totals.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$zongji.Name = "总计"

$q1.Range("B1:D1").Copy()
$zongji.Range("B1:D1").PasteSpecial(-4122)    # xlPasteFormats

$q1.Range("A2").Copy()
$zongji.Range("A2:A7").PasteSpecial(-4122)    # xlPasteFormats

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $zongji.Cells.Item(1, $col).Value = $totalHeaders[$col - 2]
}

$totalRows = @(
    @("2022-Q1", 22, 8.07),
    @("2021-Q4", 13, 8.630000000000001),
    @("2021-Q3", 16, 9.130000000000001),
    @("2021-Q2", 5, 5.01),
    @("2021-Q1", 22, 0.99),
    @("2020-Q4", 14, 0.45)
)

$row = 2
foreach ($t in $totalRows) {
    $zongji.Cells.Item($row, 1).Value = $row - 2
    $zongji.Cells.Item($row, 2).Value = $t[0]
    $zongji.Cells.Item($row, 3).Value = $t[1]
    $zongji.Cells.Item($row, 4).Value = $t[2]
    $row = $row + 1
}

$q1.Range("A1").Select() | Out-Null
$zongji.Range("A1").Select() | Out-Null
